$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "nom" column (FW).
# This shifts the old FW ("nom") -> FX and old FX ("url_produit") -> FY,
# and makes room for a new latest-snapshot price column at FW.
$ws.Columns("FW:FW").Insert()

# Populate the new FW column (rows 2-209) with a copy of the previous
# latest snapshot column (FV), which carries the same price value
# (or stays blank for rows without current price data).
$ws.Range("FV2:FV209").Copy($ws.Range("FW2:FW209"))

# Set the new header cell (row 1) to the new scrape timestamp.
$ws.Range("FW1").Value = "2026-02-05 07:42:32"
